# The deck's final slide (slide 5 / sldId 1064 - the "talentpath.com"
# title-card slide with the background picture + logo) is removed.
# Deleting it through the Slides collection also drops its associated
# notes page (notesSlide4.xml) and the now-unused Content_Types /
# relationship entries for both parts.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
